$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Lorenzo Zuani"
$ws.Range("B12").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C12").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("D12").Value = "Federico Nicolodi | U.SGUARNA"
$ws.Range("E12").Value = "Filippo Benetti | I Magnifici"
$ws.Range("F12").Value = "Nicholas Marzadro | SBARX"
